{"js": "// The author re-edited text in the \"Univariate\" intro paragraph. Word\n// automatically tracks the location of the most recent edit with a\n// hidden bookmark named \"_GoBack\" (used for Shift+F5 / \"go back to last\n// edit\"). Re-editing that paragraph therefore:\n//   1. Removes the old \"_GoBack\" bookmark from wherever it previously was\n//      (an empty paragraph near the end of the document).\n//   2. Inserts a new \"_GoBack\" bookmark around the edited text, here\n//      wrapping the first part of the sentence: \"The univariate analysis\n//      presents \".\n// (All the other bookmarks' numeric w:id values simply get renumbered in\n// document order as a side effect when Word re-saves the file; their\n// names/positions are unchanged, so nothing else needs to be touched.)\n\nconst doc = context.document;\n\n// 1. Remove the existing \"_GoBack\" bookmark, if present.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Find the target text and wrap the first clause with a new\n//    \"_GoBack\" bookmark.\nconst results = doc.body.search(\"The univariate analysis presents \", {\n  matchCase: true,\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author re-edited text in the \"Univariate\" intro paragraph. Word\n# automatically tracks the location of the most recent edit with a hidden\n# bookmark named \"_GoBack\" (used for Shift+F5 / \"go back to last edit\").\n# Re-editing that paragraph therefore moves \"_GoBack\" from its old spot\n# (an empty paragraph near the end of the document) onto the newly\n# edited text, here wrapping the first part of the sentence: \"The\n# univariate analysis presents \".\n#\n# Word's Bookmarks.Add re-seats a bookmark of the same name rather than\n# creating a duplicate, so adding \"_GoBack\" at the new range both removes\n# it from its old location and places it at the new one. All other\n# bookmarks' numeric IDs are then renumbered in document order as a\n# side effect of the save -- their names/positions don't change.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.Execute(\"The univariate analysis presents \")\n\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n"}
